$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Boca: " label - remove the stray leading space
$ws.Range("A5").Value = "Boca: "

# Move the selection to A6 (reflecting where the user left off)
$ws.Range("A6").Select()
